$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 13 and 14: Polkadot and WrappedEther swap ranking positions ---
# Row 13 becomes Polkadot (was WrappedEther), Row 14 becomes WrappedEther (was Polkadot)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.888"
$ws.Range("D13").Style = "Normal"
# E13 unchanged ("  -1.44%  ")

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.830.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

# --- Remaining price / volume updates ---

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.473.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.16"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4589"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.71%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07896"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9695"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.27%  "

$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06640"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("E19").Value = "  -1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.458.98"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.341"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.069.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.059"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.276"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9462"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.325"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05930"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.160"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.031"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5769"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1831"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5450"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.870"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("E51").Value = "  -1.35%  "
